$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "G0:Axxx:Bxxx:Zxxx"  ->  Move motor
$ws.Range("C4").Value = "G"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = ":"
$ws.Range("F4").Value = "A"
$ws.Range("G4").Value = "xxx"
$ws.Range("H4").Value = ":"
$ws.Range("I4").Value = "B"
$ws.Range("J4").Value = "xxx"
$ws.Range("K4").Value = ":"
$ws.Range("L4").Value = "Z"
$ws.Range("M4").Value = "xxx"
$ws.Range("O4").Value = "Move motor"

# Row 7: "G1:Axxx:Bxxx:Zxxx"  ->  Change speed
$ws.Range("C7").Value = "G"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = ":"
$ws.Range("F7").Value = "A"
$ws.Range("G7").Value = "xxx"
$ws.Range("H7").Value = ":"
$ws.Range("I7").Value = "B"
$ws.Range("J7").Value = "xxx"
$ws.Range("K7").Value = ":"
$ws.Range("L7").Value = "Z"
$ws.Range("M7").Value = "xxx"
$ws.Range("O7").Value = "Change speed"

# Column widths (narrow, bestFit-style columns used for the protocol string)
$ws.Columns("C:D").ColumnWidth = 1.1666666666666667
$ws.Columns("E:E").ColumnWidth = 0.5
$ws.Columns("F:F").ColumnWidth = 1.0
$ws.Columns("G:G").ColumnWidth = 2.5
$ws.Columns("H:H").ColumnWidth = 0.5
$ws.Columns("I:I").ColumnWidth = 1.0
$ws.Columns("J:J").ColumnWidth = 2.5
$ws.Columns("K:K").ColumnWidth = 0.5
$ws.Columns("L:L").ColumnWidth = 0.8333333333333334
$ws.Columns("M:M").ColumnWidth = 2.5

# Selection + zoom to match the saved view state
$ws.Range("C10").Select()
$excel.ActiveWindow.Zoom = 160
